# Fruta / hortaliza, semanal
# The underlying weekly data rows (2-15) have been reshuffled: the values in
# columns D (Fecha), M (Volumen), N (Precio minimo), O (Precio maximo),
# P (Precio promedio ponderado), Q (Unidad de comercializacion) and
# S (Precio $/Kg) for each row get replaced with the values that used to sit
# in a different row. Columns A,B,C,E,F,G,H,I,J,K,L,R,T are identical for
# every data row so they are unaffected by the shuffle.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "before" values for the columns that move, keyed by row number.
$cols = @("D", "M", "N", "O", "P", "Q", "S")
$before = @{}
for ($r = 2; $r -le 15; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $before[$r] = $rowVals
}

# Mapping of new row -> source row (where the values now placed in the new
# row used to live before the edit). Rows not listed keep their own values.
$mapping = @{
    3  = 10
    4  = 3
    5  = 9
    6  = 13
    7  = 8
    8  = 4
    9  = 5
    10 = 6
    12 = 15
    13 = 7
    15 = 12
}

foreach ($newRow in $mapping.Keys) {
    $srcRow = $mapping[$newRow]
    $srcVals = $before[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$newRow").Value = $srcVals[$c]
    }
}
